# Weekly update: insert this week's new price record at the top of the
# data table (row 79) and push the existing historical rows down by one.
#
# The sheet holds one flat table (header in row 1, data from row 2..155).
# A new "Mango" price observation is added; every existing record from the
# old row 79 onward shifts down one row (79->80, ..., 155->156), which also
# grows the used range from T155 to T156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 79:155 down to 80:156, leaving a blank row 79 to fill in.
$ws.Rows.Item(79).Insert()

# New weekly record for row 79.
$ws.Cells.Item(79, 1).Value  = 4
$ws.Cells.Item(79, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value  = "Los Lagos"
$ws.Cells.Item(79, 4).Value  = 44589
$ws.Cells.Item(79, 5).Value  = 10
$ws.Cells.Item(79, 6).Value  = "Fruta"
$ws.Cells.Item(79, 7).Value  = 100108
$ws.Cells.Item(79, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value  = 100108002
$ws.Cells.Item(79, 10).Value = "Mango"
$ws.Cells.Item(79, 11).Value = "Sin especificar"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 160
$ws.Cells.Item(79, 14).Value = 8000
$ws.Cells.Item(79, 15).Value = 8500
$ws.Cells.Item(79, 16).Value = 8250
$ws.Cells.Item(79, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(79, 18).Value = "Perú"
$ws.Cells.Item(79, 19).Value = 2062
$ws.Cells.Item(79, 20).Value = 4
